# Restore cell C10 on the active (first) worksheet from 18 to 1,
# matching revision 0c920b07a1b2dcd5a0c4cef35ff917b6ae2711fa.TEST
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("C10").Value = 1
